$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2 through 46
# from serial date 45172 (2023-09-03) to 45175 (2023-09-06)
for ($row = 2; $row -le 46; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45172) {
        $cell.Value2 = 45175
    }
}
